$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 11 (5th log entry) with the new time-tracking data
$ws.Range("B11").Value = 43867
$ws.Range("C11").Value = 0.79166666666666663
$ws.Range("D11").Value = 0.93055555555555547
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = 200
$ws.Range("G11").Value = "Kodutöö Razor pages"
$ws.Range("J11").Value = 4

# Update the active selection to match the new edit location
$ws.Range("J12").Select() | Out-Null
